# Add a new "national_id" column (M) to the Payment Plan payment list export,
# mirroring the formatting of the existing "status" column (L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column L (status) formatting/styles into the new column M for all the
# rows that currently have data/formatting (header + 2 data rows + 6 blank
# styled rows = rows 1-10).
$srcRange = $ws.Range("L1:L10")
$dstRange = $ws.Range("M1:M10")
$srcRange.Copy($dstRange)

# Match column M's width to column L's width (both render as 26.5 units wide).
$ws.Columns.Item(13).ColumnWidth = $ws.Columns.Item(12).ColumnWidth

# Header for the new column.
$ws.Range("M1").Value2 = "national_id"

# First data row (714a72db-...) has no national id value.
$ws.Range("M2").ClearContents()

# Second data row (a15e9214-...) has national id "321".
$ws.Range("M3").Value2 = "321"
